# BBI-23-6.xlsx — "Add files via upload" edit
# Fill in missing lab-4 (G) and lab-5 (H) scores for a few students, and
# recompute the selection/active-cell position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Арзамаскина Юлиана Анатольевна): Лаба №4 and Лаба №5 scores entered
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 5

# Row 7 (Гогчян Армен Давидович): Лаба №4 score entered
$ws.Range("G7").Value = 0

# Row 9 (Захаренкова Екатерина Денисовна): Лаба №4 score entered
$ws.Range("G9").Value = 0

# Row 10 (Иванов Дмитрий Сергеевич): Лаба №4 score entered
$ws.Range("G10").Value = 0

# Row 24 (Шаблыгин Михаил Максимович): Лаба №2 score entered
$ws.Range("E24").Value = 4

# Row 25 (неявка entry): Лаба №2 score entered
$ws.Range("E25").Value = 4

# Move the active cell / selection to G8, matching the edited region
$ws.Range("G8").Select() | Out-Null
